# Reverse the order of the comma-separated "Recorded By" names in column G.
# Cells with a single name (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "" -and $val.ToString().Contains(",")) {
        $parts = $val.ToString().Split(",")
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }

        $reversed = @()
        for ($i = $trimmed.Length - 1; $i -ge 0; $i--) {
            $reversed += $trimmed[$i]
        }

        $cell.Value2 = [string]::Join(", ", $reversed)
    }
}
